$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.798.52"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "2.028.64"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'227.20"
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("D6").Value = "'0.612"
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("D7").Value = "'59.76"
$ws.Range("E7").Value = "  +2.04%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("D10").Value = "'0.0814"
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").Value = "2.331.10"
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("D13").Value = "'14.56"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").Value = "'21.23"
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("D15").Value = "'0.761"
$ws.Range("E15").Value = "  +1.31%  "
$ws.Range("D16").Value = "'5.16"
$ws.Range("E16").Value = "  -2.64%  "
$ws.Range("D17").Value = "2.030.70"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").Value = "37.756.32"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").Value = "'6.00"
$ws.Range("E19").Value = "  -3.58%  "
$ws.Range("D20").Value = "'70.06"
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").Value = "0.0₃0824"
$ws.Range("E21").Value = "  -1.37%  "
$ws.Range("D22").Value = "'224.62"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "'2.42"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("D25").Value = "'2.24"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'166.13"
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'9.27"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("E28").Value = "  -3.79%  "
$ws.Range("D29").Value = "'18.89"
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("D30").Value = "'1.28"
$ws.Range("E30").Value = "  -5.22%  "
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("D32").Value = "'4.42"
$ws.Range("E32").Value = "  -2.70%  "
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0609"
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "'4.54"
$ws.Range("E35").Value = "  -1.19%  "
$ws.Range("D36").Value = "'6.41"
$ws.Range("E36").Value = "  +6.27%  "
$ws.Range("E37").Value = "  -3.35%  "
$ws.Range("D38").Value = "'3.27"
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").Value = "1.525.44"
$ws.Range("E40").Value = "  +2.51%  "
$ws.Range("D41").Value = "'0.0218"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").Value = "'96.57"
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("D43").Value = "'16.79"
$ws.Range("E43").Value = "  +1.31%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").Value = "'0.0915"
$ws.Range("E45").Value = "  -2.10%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "'1.11"
$ws.Range("E46").Value = "  -1.91%  "
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").Value = "'4.03"
$ws.Range("E47").Value = "  -1.89%  "
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("D51").Value = "2.220.06"
$ws.Range("E51").Value = "  -0.76%  "
